$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.938.24'
$ws.Range("E2").Value = '  -1.68%  '

# Row 3
$ws.Range("D3").Value = '2.294.59'
$ws.Range("E3").Value = '  -2.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.58'
$ws.Range("E5").Value = '  -6.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.90'
$ws.Range("E6").Value = '  +3.64%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -2.19%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("E9").Value = '  -4.52%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.74'
$ws.Range("E10").Value = '  -0.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").Value = '  -1.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.23'
$ws.Range("E12").Value = '  -2.55%  '

# Row 13
$ws.Range("E13").Value = '  -0.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.966'
$ws.Range("E14").Value = '  -3.58%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.32'
$ws.Range("E15").Value = '  -6.07%  '

# Row 16
$ws.Range("D16").Value = '2.636.77'
$ws.Range("E16").Value = '  -3.35%  '

# Row 17
$ws.Range("D17").Value = '2.293.69'
$ws.Range("E17").Value = '  -3.17%  '

# Row 18
$ws.Range("D18").Value = '41.925.86'
$ws.Range("E18").Value = '  -1.52%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").Value = '  -3.01%  '

# Row 20
$ws.Range("E20").Value = '  -2.99%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.05'
$ws.Range("E21").Value = '  -2.00%  '

# Row 22
$ws.Range("E22").Value = '  -8.70%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.26'
$ws.Range("E23").Value = '  -4.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  -2.90%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.16'
$ws.Range("E25").Value = '  -8.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.57%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.95'
$ws.Range("E27").Value = '  -4.61%  '

# Row 28
$ws.Range("E28").Value = '  +2.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.68'
$ws.Range("E29").Value = '  -2.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.31'
$ws.Range("E30").Value = '  -0.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.16'
$ws.Range("E31").Value = '  -6.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0881'
$ws.Range("E32").Value = '  -2.53%  '

# Row 33
$ws.Range("E33").Value = '  -6.34%  '

# Row 34
$ws.Range("E34").Value = '  -4.91%  '

# Row 35
$ws.Range("E35").Value = '  -4.15%  '

# Row 36
$ws.Range("E36").Value = '  +9.17%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.48'
$ws.Range("E37").Value = '  -3.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0349'
$ws.Range("E38").Value = '  -3.06%  '

# Row 39
$ws.Range("E39").Value = '  -5.64%  '

# Row 40
$ws.Range("E40").Value = '  -7.91%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.66'
$ws.Range("E41").Value = '  +9.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.46'
$ws.Range("E42").Value = '  -4.65%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.44'
$ws.Range("E43").Value = '  -0.75%  '

# Row 44
$ws.Range("E44").Value = '  -2.44%  '

# Row 45
$ws.Range("E45").Value = '  +0.23%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.03'
$ws.Range("E46").Value = '  +0.26%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '110.60'
$ws.Range("E47").Value = '  -5.98%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.32'
$ws.Range("E48").Value = '  -2.80%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.87'
$ws.Range("E49").Value = '  -2.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.74'
$ws.Range("E50").Value = '  +5.30%  '

# Row 51
$ws.Range("E51").Value = '  -0.87%  '
